$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 459, shifting rows 459:487 down to 460:488.
$ws.Rows(459).Insert()

# The newly inserted row 459 copies formatting from the row above (458) by default
# in Excel; explicitly set the values for the new weekly entry.
$ws.Cells.Item(459, 1).Value = 1
$ws.Cells.Item(459, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(459, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(459, 4).Value = [DateTime]"2023-01-05"
$ws.Cells.Item(459, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(459, 5).Value = 15
$ws.Cells.Item(459, 6).Value = 100112023
$ws.Cells.Item(459, 7).Value = "Brócoli"
$ws.Cells.Item(459, 8).Value = "Sin especificar"
$ws.Cells.Item(459, 9).Value = "Segunda"
$ws.Cells.Item(459, 10).Value = 1200
$ws.Cells.Item(459, 11).Value = 600
$ws.Cells.Item(459, 12).Value = 700
$ws.Cells.Item(459, 13).Value = 650
$ws.Cells.Item(459, 14).Value = "$/unidad"
$ws.Cells.Item(459, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(459, 16).Value = 650
$ws.Cells.Item(459, 17).Value = 1
$ws.Cells.Item(459, 18).Value = "Hortaliza"
